$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 10
$ws.Range("C2").Value = 100
$ws.Range("D3").Value = 100

$ws.Range("B1").Select()
